$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "57.339.50"
$ws.Cells.Item(2, 5).Value = "  -4.13%  "

$ws.Cells.Item(3, 4).Value = "2.940.47"
$ws.Cells.Item(3, 5).Value = "  -0.77%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.00%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "556.19"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -2.62%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "131.36"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +5.37%  "

$ws.Cells.Item(7, 5).Value = "  -0.14%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.513"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +2.27%  "

$ws.Cells.Item(9, 4).Value = "2.934.74"
$ws.Cells.Item(9, 5).Value = "  -0.72%  "

$ws.Cells.Item(10, 5).Value = "  -3.09%  "

$ws.Cells.Item(11, 5).Value = "  -5.66%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.445"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +1.61%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000221"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.22%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "32.50"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.39%  "

$ws.Cells.Item(15, 5).Value = "  +1.24%  "

$ws.Cells.Item(16, 4).Value = "3.422.08"
$ws.Cells.Item(16, 5).Value = "  -0.86%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "6.82"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +10.70%  "

$ws.Cells.Item(18, 4).Value = "2.931.26"
$ws.Cells.Item(18, 5).Value = "  -1.09%  "

$ws.Cells.Item(19, 4).Value = "57.373.12"
$ws.Cells.Item(19, 5).Value = "  -4.00%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "416.45"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -3.69%  "

$ws.Cells.Item(21, 5).Value = "  +0.19%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "0.681"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +3.32%  "

$ws.Cells.Item(23, 5).Value = "  -0.51%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "12.95"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.97%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "78.99"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.06%  "

$ws.Cells.Item(26, 5).Value = "  +0.03%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.03%  "

$ws.Cells.Item(28, 5).Value = "  -1.45%  "

$ws.Cells.Item(29, 5).Value = "  +4.11%  "

$ws.Cells.Item(30, 5).Value = "  +5.09%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "6.06"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.83%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "25.03"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.89%  "

$ws.Cells.Item(33, 5).Value = "  +9.58%  "

$ws.Cells.Item(34, 5).Value = "  +1.14%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.936"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.58%  "

$ws.Cells.Item(36, 5).Value = "  -3.46%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "48.44"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.27%  "

$ws.Cells.Item(38, 4).Value = "0.0₃0678"
$ws.Cells.Item(38, 5).Value = "  +4.00%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "8.40"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +6.11%  "

$ws.Cells.Item(40, 5).Value = "  +4.03%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.0348"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -2.32%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.108"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.92%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "375.88"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.97%  "

$ws.Cells.Item(44, 4).Value = "2.634.82"
$ws.Cells.Item(44, 5).Value = "  +0.58%  "

$ws.Cells.Item(45, 5).Value = "  -0.01%  "

$ws.Cells.Item(46, 5).Value = "  +1.95%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "122.26"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +3.30%  "

$ws.Cells.Item(48, 5).Value = "  +2.59%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.98"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.37%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "23.29"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.05%  "

$ws.Cells.Item(51, 5).Value = "  +0.66%  "

